{"js": "// Replace the four \"Augmentations\" bullet texts with their updated wording.\n// Each old string is unique in the document body, so body.search with\n// matchCase (and no wildcards) reliably targets exactly one run of text.\nconst replacements = [\n  {\n    find: \"Augmentations:\",\n    replace: \"Augmentations (Sets: Semiotics / DCI / DOM Statements Source):\"\n  },\n  {\n    find: \"Schema Aggregation: Type / Relationship (Role) Inference.\",\n    replace: \"Schema Aggregation: Type / Relationships (Kinds / Roles) Inference. Clustering: Unsupervised Features Learning.\"\n  },\n  {\n    find: \"Data Alignment: Attribute Inference.\",\n    replace: \"Data Alignment: Type (Feature) Attributes Value Inference. Classification: (gender, salary range: scaling).\"\n  },\n  {\n    find: \"Behavior Activation: Available Transforms (State Browsing) Inference.\",\n    replace: \"Behavior Activation: Available Transforms (State Browsing) Inference. Regression: State (class attributes values in scenario: relationship flow).\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the four \"Augmentations\" bullet paragraphs with their revised wording.\n# Each old string is unique in the document, so a plain Find/Replace (no\n# wildcards) targets exactly one paragraph per call.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Augmentations:\"; Replace = \"Augmentations (Sets: Semiotics / DCI / DOM Statements Source):\" },\n    @{ Find = \"Schema Aggregation: Type / Relationship (Role) Inference.\"; Replace = \"Schema Aggregation: Type / Relationships (Kinds / Roles) Inference. Clustering: Unsupervised Features Learning.\" },\n    @{ Find = \"Data Alignment: Attribute Inference.\"; Replace = \"Data Alignment: Type (Feature) Attributes Value Inference. Classification: (gender, salary range: scaling).\" },\n    @{ Find = \"Behavior Activation: Available Transforms (State Browsing) Inference.\"; Replace = \"Behavior Activation: Available Transforms (State Browsing) Inference. Regression: State (class attributes values in scenario: relationship flow).\" }\n)\n\nforeach ($item in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $item.Find\n    $range.Find.Replacement.Text = $item.Replace\n    $range.Find.Execute($item.Find, $false, $false, $false, $false, $false, $true, 1, $false, $item.Replace, 2)\n}\n"}
